# Updates cryptocurrency "Price" (column D) values per the Dec 26 2022 symbol-list refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new price (kept as text to match the original inline-string cell type)
$updates = @(
    @("D2", "242.99"),
    @("D3", "23.00"),
    @("D4", "5.390"),
    @("D5", "0.05908"),
    @("D6", "3.454"),
    @("D7", "6.548"),
    @("D8", "0.8106"),
    @("D9", "0.9143"),
    @("D10", "0.1416"),
    @("D11", "0.07433"),
    @("D12", "0.03280"),
    @("D13", "0.03065"),
    @("D14", "0.09339"),
    @("D15", "3.850"),
    @("D16", "0.001558"),
    @("D17", "0.04665"),
    @("D18", "0.0005931"),
    @("D19", "0.006002"),
    @("D21", "0.004916"),
    @("D22", "0.00009502"),
    @("D24", "2.148"),
    @("D40", "0.03950"),
    @("D41", "0.006188"),
    @("D42", "0.1073"),
    @("D43", "0.002531"),
    @("D44", "0.008100"),
    @("D45", "0.00005171"),
    @("D48", "0.002282")
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newValue = $pair[1]
    $cell = $ws.Range($cellRef)
    # Force text storage so values like "23.00" keep their exact formatting
    # instead of being coerced to the number 23.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    # Restore default styling so no extraneous formatting is introduced.
    $cell.Style = "Normal"
}
